# Assignment5_testing_times.xlsx edit script
# Rebuilds the sheet: renamed to "Rower 1", three timing tables (full /
# half / 1-10th length) plus a "File info" + "Ethan MacBook Pro Specs"
# side panel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Rename the sheet
# ---------------------------------------------------------------------
$ws.Name = "Rower 1"

# ---------------------------------------------------------------------
# 2. Seed cellXfs[1] (red font, no alignment) before any centered style
#    is created, so the style table comes out in the same order as the
#    target file (index 1 = font-only, index 2 = font+center).
# ---------------------------------------------------------------------
$ws.Range("F1").Font.Color = 255

# ---------------------------------------------------------------------
# 3. Write shared strings in the exact order needed to reproduce the
#    target workbook's shared-string table indices.
#    (map / pmap already live in A3 / A4 from the original file and are
#    left untouched so they keep indices 0 and 1.)
# ---------------------------------------------------------------------

# idx 2-6: relocate the "File info" panel from column G to column I
$ws.Range("I1").Value = "File info"
$ws.Range("G2").Value = ""
$ws.Range("I2").Value = "Datafile.txt"
$ws.Range("G3").Value = ""
$ws.Range("I3").Value = "3.5 million rows"
$ws.Range("G4").Value = ""
$ws.Range("I4").Value = "11 columns"
$ws.Range("G5").Value = ""
$ws.Range("I5").Value = "120 MB"
$ws.Range("G6").Value = ""

# idx 7: table 1 title
$ws.Range("A1").Value = "max_length = approx 3.5 millions rows"

# idx 8: Mac spec panel title
$ws.Range("F1").Value = "Ethan MacBook Pro Specs"

# idx 9-11: column headers (reused by all three tables)
$ws.Range("B2").Value = "Test 1(s)"
$ws.Range("C2").Value = "Test 2(s)"
$ws.Range("D2").Value = "Test 3(s)"

# idx 12-13: Mac spec panel details
$ws.Range("F2").Value = "2.2 GHz Quad-Core Intel Core i7"
$ws.Range("F3").Value = "16GB RAM"

# idx 14: table 2 title
$ws.Range("A6").Value = "half_length = approx 1.75 millions rows"

# idx 15: table 3 title
$ws.Range("A11").Value = "1/10_length = approx 350,000  rows"

# ---------------------------------------------------------------------
# 4. Fill in the remaining header rows (reuse idx 9-11 / 0-1) and data
# ---------------------------------------------------------------------
$ws.Range("B7").Value = "Test 1(s)"
$ws.Range("C7").Value = "Test 2(s)"
$ws.Range("D7").Value = "Test 3(s)"

$ws.Range("B12").Value = "Test 1(s)"
$ws.Range("C12").Value = "Test 2(s)"
$ws.Range("D12").Value = "Test 3(s)"

$ws.Range("A3").Value = "map"
$ws.Range("A4").Value = "pmap"
$ws.Range("A8").Value = "map"
$ws.Range("A9").Value = "pmap"
$ws.Range("A13").Value = "map"
$ws.Range("A14").Value = "pmap"

# Table 1 (full length) data
$ws.Range("B3").Value = 97.04
$ws.Range("C3").Value = 96.48
$ws.Range("D3").Value = 97.09
$ws.Range("B4").Value = 59.244
$ws.Range("C4").Value = 60.68
$ws.Range("D4").Value = 58.996

# Table 2 (half length) data
$ws.Range("B8").Value = 48.543
$ws.Range("C8").Value = 48.154
$ws.Range("D8").Value = 48.032
$ws.Range("B9").Value = 27.18
$ws.Range("C9").Value = 27.835
$ws.Range("D9").Value = 27.783

# Table 3 (1/10 length) data
$ws.Range("B13").Value = 8.926
$ws.Range("C13").Value = 9.23
$ws.Range("D13").Value = 9.154
$ws.Range("B14").Value = 5.708
$ws.Range("C14").Value = 5.573
$ws.Range("D14").Value = 5.578

# ---------------------------------------------------------------------
# 5. Styling: red font across the three tables + side panel (style s1),
#    then red font + centered on the three merged title rows (style s2)
# ---------------------------------------------------------------------
$ws.Range("A2:D4").Font.Color = 255
$ws.Range("F1:F3").Font.Color = 255
$ws.Range("A7:D9").Font.Color = 255
$ws.Range("A12:D14").Font.Color = 255

$ws.Range("A1:D1").Font.Color = 255
$ws.Range("A1:D1").HorizontalAlignment = -4108
$ws.Range("A6:D6").Font.Color = 255
$ws.Range("A6:D6").HorizontalAlignment = -4108
$ws.Range("A11:D11").Font.Color = 255
$ws.Range("A11:D11").HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# 6. Merge the three title rows
# ---------------------------------------------------------------------
$ws.Range("A1:D1").Merge()
$ws.Range("A6:D6").Merge()
$ws.Range("A11:D11").Merge()

# ---------------------------------------------------------------------
# 7. Column widths (account for the engine's internal +5/6 char offset
#    so the exported OOXML width matches the target workbook)
# ---------------------------------------------------------------------
$ws.Columns("F").ColumnWidth = 26.666666666666668
$ws.Columns("G").ColumnWidth = 13.498697916666666
$ws.Columns("I").ColumnWidth = 13.498697916666666

# ---------------------------------------------------------------------
# 8. Selection, matching the saved view state in the target file
# ---------------------------------------------------------------------
$ws.Range("D20").Select()
